# "F04 Froze Token Embeddings and Decoder 12"
# This notebook re-run produced new per-epoch validation-accuracy numbers
# in column B (rows 3-118, row 2/5/10/14/113 unchanged) plus a refreshed
# Python repr() (new object memory address) for the inline-string labels
# in column A for rows 102-118. The active cell in the frozen pane
# selection also moved from O13 to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated accuracy values (col B) -------------------------------------
$newValues = @{
    "B3" = 0.6875
    "B4" = 0.65625
    "B6" = 0.5625
    "B7" = 0.53125
    "B8" = 0.484375
    "B9" = 0.53125
    "B11" = 0.53125
    "B12" = 0.5
    "B13" = 0.5
    "B15" = 0.5
    "B16" = 0.5
    "B17" = 0.515625
    "B18" = 0.515625
    "B19" = 0.515625
    "B20" = 0.5
    "B21" = 0.515625
    "B22" = 0.515625
    "B23" = 0.515625
    "B24" = 0.515625
    "B25" = 0.53125
    "B26" = 0.53125
    "B27" = 0.53125
    "B28" = 0.53125
    "B29" = 0.53125
    "B30" = 0.53125
    "B31" = 0.53125
    "B32" = 0.53125
    "B33" = 0.53125
    "B34" = 0.53125
    "B35" = 0.53125
    "B36" = 0.53125
    "B37" = 0.53125
    "B38" = 0.53125
    "B39" = 0.53125
    "B40" = 0.53125
    "B41" = 0.53125
    "B42" = 0.53125
    "B43" = 0.53125
    "B44" = 0.53125
    "B45" = 0.53125
    "B46" = 0.53125
    "B47" = 0.53125
    "B48" = 0.53125
    "B49" = 0.53125
    "B50" = 0.53125
    "B51" = 0.53125
    "B52" = 0.53125
    "B53" = 0.53125
    "B54" = 0.53125
    "B55" = 0.53125
    "B56" = 0.53125
    "B57" = 0.53125
    "B58" = 0.53125
    "B59" = 0.53125
    "B60" = 0.53125
    "B61" = 0.53125
    "B62" = 0.515625
    "B63" = 0.515625
    "B64" = 0.515625
    "B65" = 0.515625
    "B66" = 0.515625
    "B67" = 0.515625
    "B68" = 0.515625
    "B69" = 0.515625
    "B70" = 0.515625
    "B71" = 0.515625
    "B72" = 0.515625
    "B73" = 0.515625
    "B74" = 0.515625
    "B75" = 0.515625
    "B76" = 0.515625
    "B77" = 0.515625
    "B78" = 0.515625
    "B79" = 0.515625
    "B80" = 0.515625
    "B81" = 0.515625
    "B82" = 0.515625
    "B83" = 0.515625
    "B84" = 0.515625
    "B85" = 0.515625
    "B86" = 0.515625
    "B87" = 0.515625
    "B88" = 0.515625
    "B89" = 0.515625
    "B90" = 0.515625
    "B91" = 0.515625
    "B92" = 0.515625
    "B93" = 0.515625
    "B94" = 0.515625
    "B95" = 0.515625
    "B96" = 0.515625
    "B97" = 0.515625
    "B98" = 0.515625
    "B99" = 0.515625
    "B100" = 0.515625
    "B101" = 0.515625
    "B102" = 0.515625
    "B103" = 0.515625
    "B104" = 0.375
    "B105" = 0.546875
    "B106" = 0.390625
    "B107" = 0.40625
    "B108" = 0.421875
    "B109" = 0.53125
    "B110" = 0.484375
    "B111" = 0.5
    "B112" = 0.671875
    "B114" = 0.34375
    "B115" = 0.484375
    "B116" = 0.421875
    "B117" = 0.453125
    "B118" = 0.3770491803278688
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

# --- Refreshed repr() text (col A, rows 102-118) --------------------------
$newRepr = "<__main__.DisplayOutputs object at 0x7f6f48504760>"
for ($r = 102; $r -le 118; $r++) {
    $ws.Range("A$r").Value = $newRepr
}

# --- Active-cell move (selection) ------------------------------------------
$ws.Range("A3").Select()
